# Generate Report for Handoff
# Updates the localization-status report with a fresh handoff generation run:
#  - "Overview" sheet: bump "Latest HO Xliff Generate Date" for the rows whose
#    handoff xliff was just (re)generated.
#  - "zh-cn" sheet: bump "Latest Handoff Datetime" for the same rows and mark
#    their "Priority" as "ht" (handoff type) now that the handoff file exists.
#  - "de-de" sheet: mark "Priority" as "ht" for the same rows as well.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Rows (1-based, matching the worksheet row numbers) affected by this handoff run.
$rows = @(7, 8, 9, 10, 12, 14)

foreach ($row in $rows) {
    $overview.Range("G$row").Value = "2016-08-21 02:19:52"
    $zhcn.Range("H$row").Value = "2016-08-21 02:19:45"
    $zhcn.Range("E$row").Value = "ht"
    $dede.Range("E$row").Value = "ht"
}
